$d = $word.ActiveDocument

# Each "<id>...</id>" marker in the document is currently split across three
# runs:
#   run 1: "<id>"        (Courier New, color 7f6000 - the XML-tag style)
#   run 2: "p053r_aN"     (plain text, color 000000)
#   run 3: "</id>"        (Courier New, color 7f6000)
#
# The edit collapses each triple into a single run, keeping run 1's
# formatting, with combined text "<id>p053r_N</id>" (the "a" in the id is
# dropped, e.g. p053r_a1 -> p053r_1). There are six such markers
# (p053r_a1 .. p053r_a6) in the document.

for ($i = 1; $i -le 6; $i++) {
    $old = "<id>p053r_a$i</id>"
    $new = "<id>p053r_$i</id>"

    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $rng.Find.Replacement.ClearFormatting()
    [void]$rng.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
